$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01_1_diccionario")

# Replace the =LOWER(A#) formulas in column B with their static, corrected values
# (the proposed names now drop the d_/f_ prefix used in column A)
$ws.Range("B2").Value = "fecha_registro"
$ws.Range("B3").Value = "ndd"
$ws.Range("B4").Value = "anio_registro"
$ws.Range("B5").Value = "fecha_incidente"
$ws.Range("B6").Value = "mes_registro"
$ws.Range("B7").Value = "mes_incidente"
$ws.Range("B8").Value = "articulo"
$ws.Range("B9").Value = "delito"
$ws.Range("B10").Value = "delito_circunstancial"
$ws.Range("B11").Value = "estado_procesal"
$ws.Range("B12").Value = "etapa_actual"
$ws.Range("B13").Value = "provincia_incidente"
$ws.Range("B14").Value = "canton_incidente"
$ws.Range("B15").Value = "tipo_delito"
$ws.Range("B16").Value = "tipo_flagrante"
$ws.Range("B17").Value = "estado_ndd"
$ws.Range("B18").Value = "grupo_horainc"

# Fix typo "Granja horaria del incidente" -> "Franja horaria del incidente"
$ws.Range("C18").Value = "Franja horaria del incidente"

# Highlight the variables still pending classification/review in yellow (10pt font)
$highlightRanges = @("A8", "A9", "A10", "A11", "A12", "A15", "A16", "A17", "A18")
foreach ($addr in $highlightRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Size = 10
    $rng.Interior.Color = 65535
}

# Update the active selection on the frozen pane to B17
$ws.Range("B17").Select()
